$d = $word.ActiveDocument

# --- Change 1: "...WordPress " + "y" (two runs) -> "...WordPress" (one run) ---
# Find/replace across the run boundary merges the two runs into a single run
# and drops the trailing space + stray "y".
$rng1 = $d.Content
$rng1.Find.Execute("Frameworks/Runtimes/CMSs: Angular, React, Node, NestJS, NextJS, Express, CodeIgniter, Drupal, WordPress y", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "Frameworks/Runtimes/CMSs: Angular, React, Node, NestJS, NextJS, Express, CodeIgniter, Drupal, WordPress", 2)

# --- Change 2: "Figma" (one run) -> "Figm" + "a" (two runs) ---
# Locate "Figma" and split it into two runs ("Figm" / "a") by toggling a
# character-level format on the trailing "a" (set then reset to its
# original value), which forces Word to break the run at that boundary
# without altering the visible formatting.
$rng2 = $d.Content
$rng2.Find.Execute("Figma", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($rng2.Find.Found) {
    $splitStart = $rng2.Start + 4
    $splitEnd = $rng2.Start + 5
    $tail = $d.Range($splitStart, $splitEnd)
    $tail.Bold = $true
    $tail.Bold = $false
}
